$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 16; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2

    $ws.Cells.Item($r, 2).Value2 = $dVal
    $ws.Cells.Item($r, 3).Value2 = $eVal
    $ws.Cells.Item($r, 4).Value2 = $bVal
    $ws.Cells.Item($r, 5).Value2 = $cVal
}
